$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Education section: add trailing spaces to two lines ---
Replace-Text "Digital Media Software Engineering" "Digital Media Software Engineering "
Replace-Text "Ferris State University • Michigan • August-2024" "Ferris State University • Michigan • August-2024 "

# --- Masetto Logistics job title: add trailing space ---
Replace-Text "Masetto Logistics" "Masetto Logistics "

# --- Masetto Logistics bullet points: full text rewrites ---
Replace-Text "Designed and developed features in Python, utilizing API calls for real-time updates for fleet management." "Developed features in Python for real-time fleet management updates, utilizing API calls."
Replace-Text "Created React front-end application for user-friendly interface and fleet visualization on an interactive map." "Created React front-end for user-friendly interface with real-time fleet location visualization."
Replace-Text "Implemented comprehensive fleet monitoring solution, enhancing operational capabilities through real-time tracking." "Implemented efficient fleet monitoring solutions for enhanced operational capabilities."
Replace-Text "Collaborated in daily standup meetings and biweekly scrum gatherings, fostering agile teamwork and communication." "Collaborated in daily standup meetings with Architecture and Front-End teams for agile teamwork."

# --- Barracuda job title: add trailing space ---
Replace-Text "Barracuda" "Barracuda "

# --- Barracuda bullet points: full text rewrites ---
Replace-Text "Created Python and Selenium automation tool to streamline error resolution for complex job executions." "Created Python and Selenium tool to automate API calls, enhancing web application functionality."
Replace-Text "Executed functional tests and UI testing for the Cloud-to-Cloud Backup web application." "Conducted functional tests and automated testing for Cloud-to-Cloud Backup web application."
Replace-Text "Leveraged JavaScript and MABL test automation tool for quality assurance and performance." "Leveraged JavaScript and MABL for robust quality assurance."
Replace-Text "Crafted comprehensive test cases and performed backend testing for data protection project's soft-delete tool." "Crafted comprehensive test cases and performed backend/UI testing for data protection project."

# --- Skills section: rewrite the 11-line list down to a 6-line list ---
Replace-Text "Front end: JavaScript, React, HTML" "Front end: JavaScript, React, CSS, HTML"
Replace-Text "Backend: Java, JDBC, REST API" "Backend: Java, Spring Boot, Rest API, SOAP"
Replace-Text "Data Flow: ETL, Data Pipelines" "Machine Learning: Python, TensorFlow, Scikit-Learn"
Replace-Text "Automation: Selenium, API Calls, Python" "Database Management: MySQL, MongoDB"
Replace-Text "Security: SSL, Data Transmission" "Agile, Jira, Git, Selenium, CI/CD, TDD"

# Remove the now-obsolete skill lines ("Database Connectivity: JDBC" through
# "Emerging Technologies: Continuous Learning", inclusive of the <w:br/>
# elements between/after them) and replace that whole span -- together with
# one of the two line breaks that used to trail the list -- with the single
# new final skill line, so only one trailing <w:br/> remains (matching the
# diff's hunk which drops a break).
$f1 = $d.Content
$f1.Find.Execute("Database Connectivity: JDBC") | Out-Null
$rangeStart = $f1.Start

$f2 = $d.Content
$f2.Find.Execute("Emerging Technologies: Continuous Learning") | Out-Null
$rangeEnd = $f2.End + 1

$toReplace = $d.Range($rangeStart, $rangeEnd)
$toReplace.Text = "AWS, SQL, NoSQL, Unit Testing, Problem-solving"
